$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl14"
$ws.Range("C2").Value = "Cxcr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.172343333333333"
$ws.Range("H2").Value = [double]"3.51703"
$ws.Range("I2").Value = [double]"0.002035925897810539"
$ws.Range("J2").Value = [double]"0.002035925897810539"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"223.2367336666667"
$ws.Range("N2").Value = [double]"669.710201"
$ws.Range("O2").Value = [double]"0.9523995969492647"
$ws.Range("P2").Value = [double]"0.9523995969492646"
$ws.Range("Q2").Value = [double]"261.7100964692255"
$ws.Range("R2").Value = [double]"2355.39086822303"
$ws.Range("S2").Value = [double]"0.001939015004493327"
$ws.Range("T2").Value = [double]"0.001939015004493327"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl14"
$ws.Range("C3").Value = "Cxcr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.172343333333333"
$ws.Range("H3").Value = [double]"3.51703"
$ws.Range("I3").Value = [double]"0.002035925897810539"
$ws.Range("J3").Value = [double]"0.002035925897810539"
$ws.Range("K3").Value = [double]"3"
$ws.Range("L3").Value = [double]"1"
$ws.Range("M3").Value = [double]"4.280784333333333"
$ws.Range("N3").Value = [double]"12.842353"
$ws.Range("O3").Value = [double]"0.01826320071400582"
$ws.Range("P3").Value = [double]"0.01826320071400582"
$ws.Range("Q3").Value = [double]"5.018548974621111"
$ws.Range("R3").Value = [double]"45.16694077159"
$ws.Range("S3").Value = [double]"3.718252331055637E-05"
$ws.Range("T3").Value = [double]"3.718252331055637E-05"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cxcl14"
$ws.Range("C4").Value = "Cxcr4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.172343333333333"
$ws.Range("H4").Value = [double]"3.51703"
$ws.Range("I4").Value = [double]"0.002035925897810539"
$ws.Range("J4").Value = [double]"0.002035925897810539"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"6.876463666666666"
$ws.Range("N4").Value = [double]"20.629391"
$ws.Range("O4").Value = [double]"0.02933720233672951"
$ws.Range("P4").Value = [double]"0.0293372023367295"
$ws.Range("Q4").Value = [double]"8.061576336525555"
$ws.Range("R4").Value = [double]"72.55418702873"
$ws.Range("S4").Value = [double]"5.972837000665546E-05"
$ws.Range("T4").Value = [double]"5.972837000665546E-05"

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl14"
$ws.Range("C5").Value = "Cxcr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"574.5849303333333"
$ws.Range("H5").Value = [double]"1723.754791"
$ws.Range("I5").Value = [double]"0.9978410819560519"
$ws.Range("J5").Value = [double]"0.9978410819560519"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"223.2367336666667"
$ws.Range("N5").Value = [double]"669.710201"
$ws.Range("O5").Value = [double]"0.9523995969492647"
$ws.Range("P5").Value = [double]"0.9523995969492646"
$ws.Range("Q5").Value = [double]"128268.4630617025"
$ws.Range("R5").Value = [double]"1154416.167555323"
$ws.Range("S5").Value = [double]"0.950343444274362"
$ws.Range("T5").Value = [double]"0.9503434442743619"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cxcl14"
$ws.Range("C6").Value = "Cxcr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"574.5849303333333"
$ws.Range("H6").Value = [double]"1723.754791"
$ws.Range("I6").Value = [double]"0.9978410819560519"
$ws.Range("J6").Value = [double]"0.9978410819560519"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"4.280784333333333"
$ws.Range("N6").Value = [double]"12.842353"
$ws.Range("O6").Value = [double]"0.01826320071400582"
$ws.Range("P6").Value = [double]"0.01826320071400582"
$ws.Range("Q6").Value = [double]"2459.674167940358"
$ws.Range("R6").Value = [double]"22137.06751146322"
$ws.Range("S6").Value = [double]"0.0182237719604441"
$ws.Range("T6").Value = [double]"0.0182237719604441"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cxcl14"
$ws.Range("C7").Value = "Cxcr4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"574.5849303333333"
$ws.Range("H7").Value = [double]"1723.754791"
$ws.Range("I7").Value = [double]"0.9978410819560519"
$ws.Range("J7").Value = [double]"0.9978410819560519"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"6.876463666666666"
$ws.Range("N7").Value = [double]"20.629391"
$ws.Range("O7").Value = [double]"0.02933720233672951"
$ws.Range("P7").Value = [double]"0.0293372023367295"
$ws.Range("Q7").Value = [double]"3951.112396851364"
$ws.Range("R7").Value = [double]"35560.01157166228"
$ws.Range("S7").Value = [double]"0.02927386572124578"
$ws.Range("T7").Value = [double]"0.02927386572124578"

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl14"
$ws.Range("C8").Value = "Cxcr4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = [double]"1"
$ws.Range("F8").Value = [double]"0.3333333333333333"
$ws.Range("G8").Value = [double]"0.07082233333333333"
$ws.Range("H8").Value = [double]"0.212467"
$ws.Range("I8").Value = [double]"0.0001229921461375399"
$ws.Range("J8").Value = [double]"0.0001229921461375399"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"223.2367336666667"
$ws.Range("N8").Value = [double]"669.710201"
$ws.Range("O8").Value = [double]"0.9523995969492647"
$ws.Range("P8").Value = [double]"0.9523995969492646"
$ws.Range("Q8").Value = [double]"15.81014636398522"
$ws.Range("R8").Value = [double]"142.291317275867"
$ws.Range("S8").Value = [double]"0.0001171376704093181"
$ws.Range("T8").Value = [double]"0.000117137670409318"

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl14"
$ws.Range("C9").Value = "Cxcr4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = [double]"1"
$ws.Range("F9").Value = [double]"0.3333333333333333"
$ws.Range("G9").Value = [double]"0.07082233333333333"
$ws.Range("H9").Value = [double]"0.212467"
$ws.Range("I9").Value = [double]"0.0001229921461375399"
$ws.Range("J9").Value = [double]"0.0001229921461375399"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"4.280784333333333"
$ws.Range("N9").Value = [double]"12.842353"
$ws.Range("O9").Value = [double]"0.01826320071400582"
$ws.Range("P9").Value = [double]"0.01826320071400582"
$ws.Range("Q9").Value = [double]"0.3031751349834444"
$ws.Range("R9").Value = [double]"2.728576214851"
$ws.Range("S9").Value = [double]"2.246230251156226E-06"
$ws.Range("T9").Value = [double]"2.246230251156226E-06"

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cxcl14"
$ws.Range("C10").Value = "Cxcr4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = [double]"1"
$ws.Range("F10").Value = [double]"0.3333333333333333"
$ws.Range("G10").Value = [double]"0.07082233333333333"
$ws.Range("H10").Value = [double]"0.212467"
$ws.Range("I10").Value = [double]"0.0001229921461375399"
$ws.Range("J10").Value = [double]"0.0001229921461375399"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"6.876463666666666"
$ws.Range("N10").Value = [double]"20.629391"
$ws.Range("O10").Value = [double]"0.02933720233672951"
$ws.Range("P10").Value = [double]"0.0293372023367295"
$ws.Range("Q10").Value = [double]"0.4870072019552222"
$ws.Range("R10").Value = [double]"4.383064817596999"
$ws.Range("S10").Value = [double]"3.608245477065612E-06"
$ws.Range("T10").Value = [double]"3.608245477065611E-06"
